$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1450048780487805
$ws.Range("V2").Value = 0.0002448603057459146
$ws.Range("Z2").Value = -0.1512229460108298
$ws.Range("AB2").Value = -617.5886514155957
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -617.5886514155957

# Row 3
$ws.Range("T3").Value = 0.1492487804878049
$ws.Range("V3").Value = 0.0001488973818309612
$ws.Range("Z3").Value = -0.1712036912957602
$ws.Range("AB3").Value = -1149.809950923937
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -1149.809950923937

# Row 4
$ws.Range("T4").Value = 0.1469268292682927
$ws.Range("V4").Value = 0.0002222807942365138
$ws.Range("Z4").Value = -0.1118592356004134
$ws.Range("AB4").Value = -503.2339207920574
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -503.2339207920574

# Row 5
$ws.Range("T5").Value = 0.1418926829268293
$ws.Range("V5").Value = 0.0002529432437181515
$ws.Range("Z5").Value = -0.1118840187357806
$ws.Range("AB5").Value = -442.3285520148156
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -442.3285520148156

# Row 6
$ws.Range("T6").Value = 0.1446439024390244
$ws.Range("V6").Value = 0.0001851607801792304
$ws.Range("Z6").Value = -0.1491028544390983
$ws.Range("AB6").Value = -805.2615369992011
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -805.2615369992011

# Row 7
$ws.Range("T7").Value = 0.1429268292682927
$ws.Range("V7").Value = 0.0003232296608680373
$ws.Range("Z7").Value = -0.1552406173418066
$ws.Range("AB7").Value = -480.2796158152874
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -480.2796158152874

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = -0.0008470748809544388
$ws.Range("AB8").Value = "-Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "-Inf"

# Row 9
$ws.Range("T9").Value = 0.1450048780487805
$ws.Range("V9").Value = 0.0002448603057459146
$ws.Range("Z9").Value = -0.06097875807370545
$ws.Range("AB9").Value = -249.0348849640888
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = -249.0348849640888

# Row 10
$ws.Range("T10").Value = 0.1492487804878049
$ws.Range("V10").Value = 0.0001488973818309612
$ws.Range("Z10").Value = -0.05501833662875388
$ws.Range("AB10").Value = -369.5050641737582
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = -369.5050641737582

# Row 11
$ws.Range("T11").Value = 0.1469268292682927
$ws.Range("V11").Value = 0.0002222807942365138
$ws.Range("Z11").Value = -0.02689165068177779
$ws.Range("AB11").Value = -120.9805407351758
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = -120.9805407351758

# Row 12
$ws.Range("T12").Value = 0.1418926829268293
$ws.Range("V12").Value = 0.0002529432437181515
$ws.Range("Z12").Value = -0.06996458571528331
$ws.Range("AB12").Value = -276.6019154606997
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = -276.6019154606997

# Row 13
$ws.Range("T13").Value = 0.1446439024390244
$ws.Range("V13").Value = 0.0001851607801792304
$ws.Range("Z13").Value = -0.02425607325827775
$ws.Range("AB13").Value = -131.0000597037805
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = -131.0000597037805

# Row 14
$ws.Range("T14").Value = 0.1429268292682927
$ws.Range("V14").Value = 0.0003232296608680373
$ws.Range("Z14").Value = -0.1008702799687271
$ws.Range("AB14").Value = -312.0699990769369
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = -312.0699990769369

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = [double]"4.010042052983081E-05"
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"

Write-Output "done"
